$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 730781075708.5688
$ws.Range("C3").Value = 245257091704.5142
$ws.Range("C4").Value = 37558555657.86328
$ws.Range("C5").Value = 33331957411.70317
$ws.Range("C6").Value = 25311182965.7804
$ws.Range("C7").Value = 13828792225.42868
$ws.Range("C8").Value = 11456884589.16068
$ws.Range("C9").Value = 9235860793.291689
$ws.Range("C10").Value = 8453209150.625892
$ws.Range("C11").Value = 8111096848.791776
$ws.Range("C12").Value = 7869724955.00743
$ws.Range("C13").Value = 7842322201.283567
$ws.Range("C14").Value = 6816510736.856397
$ws.Range("C15").Value = 6082869313.002133
$ws.Range("C16").Value = 5208165229.352201
$ws.Range("C17").Value = 5077501592.809649
$ws.Range("C18").Value = 4482371364.744591
$ws.Range("C19").Value = 3661721163.844837
$ws.Range("C20").Value = 3520144110.326622
$ws.Range("C21").Value = 3482261875.863487
$ws.Range("C22").Value = 3395851447.907824
$ws.Range("C23").Value = 3021454679.573148
$ws.Range("C24").Value = 2940131061.636749
$ws.Range("C25").Value = 2851678189.604074
$ws.Range("C26").Value = 2446990549.970589
